$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Period:" row (from - to range), previously used Joda DateTime#toString,
# now uses dateTool.format(...)
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'

# "Time" column template cell (event.serverTime), previously built a new
# org.joda.time.DateTime and called toString, now uses dateTool.format(...)
$ws.Range("A9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", event.serverTime, locale, timezone)}'
